$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Нет."
$ws.Range("C12").Value = "1. Нет.  `n2. Нет."
$ws.Range("C24").Value = "1. Нет.  `n2. Нет."
$ws.Range("C28").Value = "1. Нет.  `n2. Да."
$ws.Range("C32").Value = "1. Нет.  `n2. Нет."
$ws.Range("C34").Value = "1. Нет.  `n2. Нет."
$ws.Range("C35").Value = "1. Нет.`n2. Нет."
$ws.Range("C49").Value = "1. Нет.`n2. Нет."
$ws.Range("C50").Value = "Нет."
$ws.Range("C60").Value = "Нет."
$ws.Range("C63").Value = "1. Да`n2. Да"
$ws.Range("C89").Value = "1. Нет.  `n2. Нет."
$ws.Range("C99").Value = "1. Нет.  `n2. Нет."
$ws.Range("C116").Value = "1. Нет.  `n2. Нет."
$ws.Range("C121").Value = "1. Нет.  `n2. Нет."
$ws.Range("C129").Value = "1. Нет.  `n2. Нет."
$ws.Range("C133").Value = "1. Да  `n2. Да"
$ws.Range("C137").Value = "1. Нет.`n2. Нет."
$ws.Range("C143").Value = "1. Нет.  `n2. Нет."
$ws.Range("C152").Value = "1. Нет.`n2. Да."
$ws.Range("C154").Value = "1. Нет.`n2. Нет."
$ws.Range("C171").Value = "1. Нет.`n2. Нет."
$ws.Range("C172").Value = "Нет."
$ws.Range("C177").Value = "Нет."
$ws.Range("C188").Value = "1. Нет.  `n2. Нет."
$ws.Range("C191").Value = "Нет."
$ws.Range("C195").Value = "1. Нет.  `n2. Нет."
$ws.Range("C198").Value = "1. Нет.`n2. Нет."
$ws.Range("C204").Value = "1. Нет.`n2. Нет."
$ws.Range("C211").Value = "1. Нет.`n2. Нет."
$ws.Range("C222").Value = "1. Нет.`n2. Нет."
$ws.Range("C225").Value = "1. Нет.  `n2. Нет."
$ws.Range("C229").Value = "1. Нет.  `n2. Нет."
$ws.Range("C230").Value = "1. Нет.`n2. Нет."
$ws.Range("C235").Value = "1. Нет.  `n2. Нет."
$ws.Range("C250").Value = "1. Нет.  `n2. Нет."
$ws.Range("C259").Value = "1. Нет.  `n2. Нет."
$ws.Range("C271").Value = "1. Нет.  `n2. Нет."
$ws.Range("C278").Value = "1. Нет.`n2. Нет."
$ws.Range("C279").Value = "Нет."
$ws.Range("C288").Value = "1. Нет.  `n2. Нет."
$ws.Range("C297").Value = "Нет."
$ws.Range("C302").Value = "1. Нет.  `n2. Нет."
$ws.Range("C303").Value = "Нет."
$ws.Range("C306").Value = "Нет."
$ws.Range("C307").Value = "1. Нет.`n2. Нет."
$ws.Range("C334").Value = "1. Нет.  `n2. Нет."
$ws.Range("C339").Value = "Нет."
$ws.Range("C357").Value = "1. Нет.`n2. Нет."
$ws.Range("C359").Value = "Нет."
$ws.Range("C363").Value = "1. Нет.`n2. Нет."
$ws.Range("C377").Value = "1. Нет.  `n2. Нет."
$ws.Range("C386").Value = "Нет."
$ws.Range("C400").Value = "Пользователь: "
$ws.Range("C401").Value = "1. Нет.  `n2. Нет."
$ws.Range("C403").Value = "Нет."
$ws.Range("C408").Value = "Нет."
$ws.Range("C413").Value = "1. Нет.  `n2. Нет."
$ws.Range("C414").Value = "1. Нет.  `n2. Нет."
$ws.Range("C419").Value = "1. Нет.  `n2. Нет."
$ws.Range("C429").Value = "Нет."
$ws.Range("C430").Value = "1. Нет.`n2. Нет."
$ws.Range("C431").Value = "Нет."
$ws.Range("C434").Value = "1. Нет.`n2. Нет."
$ws.Range("C441").Value = "1. Нет.`n2. Нет."
$ws.Range("C443").Value = "1. Нет.  `n2. Нет."
$ws.Range("C449").Value = "Нет."
$ws.Range("C456").Value = "Нет."
$ws.Range("C458").Value = "1. Нет.`n2. Нет."
$ws.Range("C471").Value = "1. Нет.`n2. Нет."
$ws.Range("C474").Value = "1. Нет.  `n2. Нет."
$ws.Range("C476").Value = "1. Нет.  `n2. Нет."
$ws.Range("C481").Value = "Нет."
$ws.Range("C483").Value = "1. Нет.`n2. Нет."
$ws.Range("C488").Value = "1. Нет.`n2. Нет."
$ws.Range("C501").Value = "1. Нет.  `n2. Нет."
$ws.Range("C502").Value = "Нет."
$ws.Range("C518").Value = "1. Нет.`n2. Нет."
$ws.Range("C530").Value = "1. Нет.`n2. Нет."
$ws.Range("C531").Value = "1. Нет.  `n2. Нет."
$ws.Range("C534").Value = "1. Нет.  `n2. Нет."
$ws.Range("C536").Value = "1. Нет.  `n2. Нет."
$ws.Range("C538").Value = "1. Нет.  `n2. Нет."
$ws.Range("C539").Value = "Нет."
$ws.Range("C545").Value = "1. Нет.  `n2. Нет."
$ws.Range("C554").Value = "Пользователь: Здравствуйте! Я хотел бы узнать о статусе моего заказа.`n Оператор: Здравствуйте! Ваш заказ находится в обработке и будет отправлен в течение 3-5 рабочих дней.`n Пользователь: Это слишком долго! Я ожидал, что он будет отправлен быстрее.`n Оператор: Мы приносим извинения за задержку. Мы стараемся обработать заказы как можно быстрее.`n Пользователь: Но я ждал уже неделю! Почему так долго?`n Оператор: К сожалению, у нас возникли некоторые задержки из-за высокого спроса. Мы делаем все возможное, чтобы ускорить процесс.`n Пользователь: Это не оправдание! Я не доволен вашим обслуживанием.`n1. Да`n2. Да"
$ws.Range("C555").Value = "1. Нет.  `n2. Нет."
$ws.Range("C557").Value = "1. Нет.`n2. Нет."
$ws.Range("C560").Value = "Нет."
$ws.Range("C565").Value = "1. Нет.  `n2. Нет."
$ws.Range("C582").Value = "Нет."
$ws.Range("C586").Value = "1. Нет.  `n2. Нет."
$ws.Range("C590").Value = "1. Нет.  `n2. Нет."
$ws.Range("C591").Value = "Нет."
$ws.Range("C601").Value = "1. Нет.  `n2. Нет."
$ws.Range("C611").Value = "1. Да.  `n2. Да."
$ws.Range("C614").Value = "1. Нет.`n2. Нет."
$ws.Range("C615").Value = "Нет."
$ws.Range("C631").Value = "Нет."
$ws.Range("C653").Value = "Нет."
